$wb = $excel.ActiveWorkbook

# The RPEpUACE source values feeding the BLAPE workbook were refreshed.
# "data from RPEpUACE" holds CH4 GWP (B11) and N2O GWP (B12) used by the
# BLAPE sheet's CH4/N2O emission rows (which multiply these by B$2 via
# formulas), so updating these two cells ripples through the whole sheet.
$wsData = $wb.Worksheets.Item("data from RPEpUACE")
$wsData.Activate()
$wsData.Range("B11").Value = 0.0043636084384378205
$wsData.Range("B12").Value = 0.0001335357177073511
$wsData.Range("G11").Select()

$wsBlape = $wb.Worksheets.Item("BLAPE")
$wsBlape.Activate()
$wsBlape.Range("D17").Select()
